$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Add a new "Group contributions" section at the end of the document,
# right before the closing section properties, consisting of:
#   1) a bold+italic heading "Group contributions"
#      (paragraph mark itself is bold only, matching how Word records
#       a bold paragraph whose run text was further italicized)
#   2) a plain paragraph "Jason Allen:"
#   3) an indented paragraph (split over 4 runs) describing Jason's
#      contribution
#
# All three paragraphs are first inserted as plain text (so nothing
# inherits stray character formatting from a "current typing" state),
# then formatting is applied back onto them one at a time.
# ------------------------------------------------------------------

$lastCountBefore = $d.Paragraphs.Count

# --- Insert the three new paragraphs, plain text only ---------------
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "Group contributions"

$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "Jason Allen:"

$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "Completed the parsing and the initial (but slow) simulator. Created the search algorithm that kind of failed, but performed research and found a good machine online and helped to modify it to make it better."

# --- Paragraph 1: "Group contributions" heading formatting -----------
$headingPara = $d.Paragraphs.Item($lastCountBefore + 1)
$headingPara.Range.Font.Bold = 1
# Italicize only the visible text, not the trailing paragraph mark,
# so the paragraph mark's run properties stay bold-only.
$headingTextRange = $d.Range($headingPara.Range.Start, $headingPara.Range.End - 1)
$headingTextRange.Font.Italic = 1

# --- Paragraph 3: indent + split the contribution text into runs -----
$contribPara = $d.Paragraphs.Item($lastCountBefore + 3)
$contribPara.LeftIndent = 36

# Replace the plain text with the same text typed in four pieces,
# separated momentarily by bookmarks so the runs do not get merged
# back into a single run on save.
$cr = $contribPara.Range
$cr.Text = ""
$cr.InsertAfter("Completed the parsing and ")
$cr.Collapse(0)
$d.Bookmarks.Add("zzTmpSplit1", $cr)

$cr.InsertAfter("the ")
$cr.Collapse(0)
$d.Bookmarks.Add("zzTmpSplit2", $cr)

$cr.InsertAfter("init")
$cr.Collapse(0)
$d.Bookmarks.Add("zzTmpSplit3", $cr)

$cr.InsertAfter("ial (but slow) simulator. Created the search algorithm that kind of failed, but performed research and found a good machine online and helped to modify it to make it better.")

# Drop the scaffolding bookmarks used only to keep the runs distinct.
$d.Bookmarks("zzTmpSplit1").Delete()
$d.Bookmarks("zzTmpSplit2").Delete()
$d.Bookmarks("zzTmpSplit3").Delete()
